$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value to a cell, preserving its original
# (General) style/number-format so numeric-looking strings like "99.40"
# are not silently coerced into numbers (which would drop trailing zeros).
function Set-TextValue($addr, $value) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue "D2" "42.843.27"
Set-TextValue "E2" "  +0.11%  "
Set-TextValue "D3" "2.570.04"
Set-TextValue "E3" "  +1.58%  "
Set-TextValue "E4" "  +0.05%  "
Set-TextValue "D5" "313.01"
Set-TextValue "E5" "  -0.80%  "
Set-TextValue "D6" "99.40"
Set-TextValue "E6" "  +3.69%  "
Set-TextValue "E7" "  -0.20%  "
Set-TextValue "E8" "  -0.05%  "
Set-TextValue "E9" "  +0.16%  "
Set-TextValue "D10" "35.77"
Set-TextValue "E10" "  -0.90%  "
Set-TextValue "D12" "7.46"
Set-TextValue "E12" "  -1.26%  "
Set-TextValue "D13" "2.962.54"
Set-TextValue "E13" "  +1.51%  "
Set-TextValue "E14" "  -1.31%  "
Set-TextValue "D15" "15.92"
Set-TextValue "E15" "  +4.74%  "
Set-TextValue "D16" "2.553.15"
Set-TextValue "E16" "  -0.24%  "
Set-TextValue "E17" "  -0.83%  "
Set-TextValue "D18" "42.905.44"
Set-TextValue "E18" "  +0.13%  "
Set-TextValue "D19" "6.77"
Set-TextValue "E19" "  -1.12%  "
Set-TextValue "D20" "12.55"
Set-TextValue "E20" "  -3.05%  "
Set-TextValue "E21" "  -0.14%  "
Set-TextValue "D22" "69.58"
Set-TextValue "E22" "  -0.50%  "
Set-TextValue "D23" "249.64"
Set-TextValue "E23" "  -1.56%  "
Set-TextValue "D24" "2.95"
Set-TextValue "E24" "  +0.26%  "
Set-TextValue "E25" "  -0.07%  "
Set-TextValue "D26" "27.11"
Set-TextValue "E26" "  +1.61%  "
Set-TextValue "D27" "1.00"
Set-TextValue "E27" "  +0.03%  "
Set-TextValue "D28" "2.40"
Set-TextValue "E28" "  -1.09%  "
Set-TextValue "D29" "39.90"
Set-TextValue "E29" "  -0.95%  "
Set-TextValue "D30" "10.26"
Set-TextValue "E30" "  -1.39%  "
Set-TextValue "B31" "Filecoin"
Set-TextValue "C31" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D31" "5.81"
Set-TextValue "E31" "  -1.83%  "
Set-TextValue "B32" "Monero"
Set-TextValue "C32" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D32" "156.87"
Set-TextValue "E32" "  -1.04%  "
Set-TextValue "E33" "  +0.97%  "
Set-TextValue "D34" "0.0802"
Set-TextValue "E34" "  +2.72%  "
Set-TextValue "D35" "2.12"
Set-TextValue "E35" "  -2.18%  "
Set-TextValue "E36" "  +0.04%  "
Set-TextValue "D37" "18.59"
Set-TextValue "E37" "  -3.14%  "
Set-TextValue "E38" "  +11.46%  "
Set-TextValue "E39" "  +0.01%  "
Set-TextValue "E40" "  -0.15%  "
Set-TextValue "D41" "23.26"
Set-TextValue "E41" "  +0.38%  "
Set-TextValue "D42" "4.11"
Set-TextValue "E42" "  +7.26%  "
Set-TextValue "B43" "VeChain"
Set-TextValue "C43" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue "D43" "0.0303"
Set-TextValue "E43" "  -0.32%  "
Set-TextValue "B44" "FirstDigitalUSD"
Set-TextValue "C44" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D44" "1.00"
Set-TextValue "E44" "  -0.06%  "
Set-TextValue "E45" "  -2.09%  "
Set-TextValue "D46" "2.011.24"
Set-TextValue "E46" "  -0.91%  "
Set-TextValue "E47" "  -1.79%  "
Set-TextValue "D48" "2.814.32"
Set-TextValue "E48" "  +1.52%  "
Set-TextValue "D49" "0.197"
Set-TextValue "E49" "  +2.83%  "
Set-TextValue "D50" "81.98"
Set-TextValue "E50" "  -3.49%  "
Set-TextValue "E51" "  -0.40%  "
